$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: updated Price (D) and Volume(1h) (E) columns.
# Cells whose new D-column text would otherwise be auto-coerced to a number
# by Excel are forced to Text format first so they stay literal strings.

$ws.Range("D2").Value = "61.378.31"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.929.08"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.84"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.26"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.97"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000224"
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.49"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "3.414.14"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "61.370.13"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "2.931.25"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "431.36"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.47"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.07"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.84"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.86"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.76"
$ws.Range("E26").Value = "  -2.31%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  -5.25%  "
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.60"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "0.0₃0884"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.42"
$ws.Range("E41").Value = "  +6.88%  "
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0346"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "2.698.20"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "133.56"
$ws.Range("E45").Value = "  +2.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "362.75"
$ws.Range("E46").Value = "  -3.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.60"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.123"
$ws.Range("E51").Value = "  -1.94%  "
